$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current values of rows 2-25, columns A-F, before reordering
$original = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 6; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value2)
    }
    $original[$r] = $rowVals
}

# Mapping of target row -> source row (the row reorder described by the fix)
$mapping = [ordered]@{
    2 = 8
    3 = 11
    4 = 4
    5 = 12
    6 = 5
    7 = 10
    8 = 6
    9 = 15
    10 = 9
    11 = 2
    12 = 3
    13 = 7
    14 = 14
    15 = 13
    16 = 19
    17 = 16
    18 = 20
    19 = 17
    20 = 18
    21 = 21
    22 = 23
    23 = 22
    24 = 24
    25 = 25
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $vals = $original[$source]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($target, $c).Value = $vals[$c - 1]
    }
}
